# Update AVERAGE_10_9 qoq forecast error table with refreshed
# ifo GDP component analysis preprocessing values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.03133494016881249
$ws.Range("C2").Value = 0.618674247399176
$ws.Range("D2").Value = 0.7992451844515852
$ws.Range("E2").Value = 0.8940051367031318
$ws.Range("F2").Value = 0.9025266721007571
$ws.Range("G2").Value = 50

$ws.Range("B3").Value = 0.05905192170129663
$ws.Range("C3").Value = 0.7500407276805185
$ws.Range("D3").Value = 1.336363517888593
$ws.Range("E3").Value = 1.156011902139677
$ws.Range("F3").Value = 1.166466736141374
$ws.Range("G3").Value = 49

$ws.Range("B4").Value = 0.003554835396229139
$ws.Range("C4").Value = 0.6262718728555511
$ws.Range("D4").Value = 0.7969898074296445
$ws.Range("E4").Value = 0.8927428562747756
$ws.Range("F4").Value = 0.9021829812466056
$ws.Range("G4").Value = 48

$ws.Range("B5").Value = -0.03511572363863584
$ws.Range("C5").Value = 0.6635828519735995
$ws.Range("D5").Value = 0.8873739543468248
$ws.Range("E5").Value = 0.9420052836087623
$ws.Range("F5").Value = 0.9517525198146293

$ws.Range("B6").Value = -0.02357502177020739
$ws.Range("C6").Value = 0.6102703326053934
$ws.Range("D6").Value = 0.7645395854549295
$ws.Range("E6").Value = 0.8743795431361198
$ws.Range("F6").Value = 0.8837201036899427
$ws.Range("G6").Value = 46

$ws.Range("B7").Value = 0.01163558867306234
$ws.Range("C7").Value = 0.5783626329962743
$ws.Range("D7").Value = 0.6903657101012227
$ws.Range("E7").Value = 0.8308824887438793
$ws.Range("F7").Value = 0.8432949623807621
$ws.Range("G7").Value = 34

$ws.Range("B8").Value = 0.009654851496814876
$ws.Range("C8").Value = 0.5447572585674562
$ws.Range("D8").Value = 0.6494921422519612
$ws.Range("E8").Value = 0.8059107532797668
$ws.Range("F8").Value = 0.8183475071966614
$ws.Range("G8").Value = 33

$ws.Range("B9").Value = -0.04525758420911288
$ws.Range("C9").Value = 0.5348296361048615
$ws.Range("D9").Value = 0.5246369517211453
$ws.Range("E9").Value = 0.7243182668697133
$ws.Range("F9").Value = 0.7466109761084062
$ws.Range("G9").Value = 16

$ws.Range("B10").Value = -0.120014643191906
$ws.Range("C10").Value = 0.3396452386262637
$ws.Range("D10").Value = 0.2564608965434484
$ws.Range("E10").Value = 0.5064196841982432
$ws.Range("F10").Value = 0.5186064467643261
$ws.Range("G10").Value = 10
